$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9161549210548401
$ws.Range("B1").Value = 2.573136806488037
$ws.Range("C1").Value = 2.65956974029541
$ws.Range("D1").Value = 2.635978937149048
$ws.Range("E1").Value = 1.891165852546692
